$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update values
$ws.Range("D2").Value = 882
$ws.Range("E2").Value = 69
$ws.Range("F2").Value = 69
$ws.Range("G2").Value = 81
$ws.Range("H2").Value = 65
$ws.Range("I2").Value = 65
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 966
$ws.Range("L2").Value = 311
$ws.Range("M2").Value = 655
$ws.Range("N2").Value = 655
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 50
$ws.Range("Q2").Value = 101
$ws.Range("R2").Value = -75
$ws.Range("S2").Value = 7
$ws.Range("T2").Value = 41
$ws.Range("U2").Value = 59
$ws.Range("V2").Value = 178
$ws.Range("W2").Value = 7.78
$ws.Range("X2").Value = 7.38
$ws.Range("Y2").Value = 10.19
$ws.Range("Z2").Value = 6.97
$ws.Range("AA2").Value = 47.47
$ws.Range("AB2").Value = 1560.11
$ws.Range("AC2").Value = 650
$ws.Range("AD2").Value = 12.73
$ws.Range("AE2").Value = 8733
$ws.Range("AF2").Value = 0.95
$ws.Range("AG2").Value = 200
$ws.Range("AH2").Value = 2.42
$ws.Range("AI2").Value = 23.07
$ws.Range("AJ2").Value = 10000000

# Row 3: update values
$ws.Range("D3").Value = 866
$ws.Range("E3").Value = 41
$ws.Range("F3").Value = 41
$ws.Range("G3").Value = 33
$ws.Range("H3").Value = 26
$ws.Range("I3").Value = 26
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 974
$ws.Range("L3").Value = 306
$ws.Range("M3").Value = 668
$ws.Range("N3").Value = 667
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 50
$ws.Range("Q3").Value = 88
$ws.Range("R3").Value = -42
$ws.Range("S3").Value = -1
$ws.Range("T3").Value = 25
$ws.Range("U3").Value = 63
$ws.Range("V3").Value = 197
$ws.Range("W3").Value = 4.77
$ws.Range("X3").Value = 3.04
$ws.Range("Y3").Value = 4
$ws.Range("Z3").Value = 2.72
$ws.Range("AA3").Value = 45.77
$ws.Range("AB3").Value = 1586.75
$ws.Range("AC3").Value = 264
$ws.Range("AD3").Value = 30.82
$ws.Range("AE3").Value = 8896
$ws.Range("AF3").Value = 0.92
$ws.Range("AG3").Value = 160
$ws.Range("AH3").Value = 1.97
$ws.Range("AI3").Value = 45.43
$ws.Range("AJ3").Value = 10000000

# Row 4: update values
$ws.Range("D4").Value = 864
$ws.Range("E4").Value = 42
$ws.Range("F4").Value = 42
$ws.Range("G4").Value = 56
$ws.Range("H4").Value = 47
$ws.Range("I4").Value = 47
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 984
$ws.Range("L4").Value = 280
$ws.Range("M4").Value = 704
$ws.Range("N4").Value = 703
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 50
$ws.Range("Q4").Value = 121
$ws.Range("R4").Value = -6
$ws.Range("S4").Value = -48
$ws.Range("T4").Value = 25
$ws.Range("U4").Value = 95
$ws.Range("V4").Value = 162
$ws.Range("W4").Value = 4.83
$ws.Range("X4").Value = 5.43
$ws.Range("Y4").Value = 6.8
$ws.Range("Z4").Value = 4.8
$ws.Range("AA4").Value = 39.72
$ws.Range("AB4").Value = 1654.91
$ws.Range("AC4").Value = 466
$ws.Range("AD4").Value = 13.37
$ws.Range("AE4").Value = 9370
$ws.Range("AF4").Value = 0.66
$ws.Range("AG4").Value = 160
$ws.Range("AH4").Value = 2.57
$ws.Range("AI4").Value = 25.74
$ws.Range("AJ4").Value = 10000000

# Row 5: update values
$ws.Range("D5").Value = 825
$ws.Range("E5").Value = 63
$ws.Range("F5").Value = 63
$ws.Range("G5").Value = 53
$ws.Range("H5").Value = 40
$ws.Range("I5").Value = 40
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 990
$ws.Range("L5").Value = 262
$ws.Range("M5").Value = 728
$ws.Range("N5").Value = 726
$ws.Range("O5").Value = 2
$ws.Range("P5").Value = 50
$ws.Range("Q5").Value = 62
$ws.Range("R5").Value = -81
$ws.Range("S5").Value = -27
$ws.Range("T5").Value = 41
$ws.Range("U5").Value = 20
$ws.Range("V5").Value = 134
$ws.Range("W5").Value = 7.63
$ws.Range("X5").Value = 4.87
$ws.Range("Y5").Value = 5.6
$ws.Range("Z5").Value = 4.07
$ws.Range("AA5").Value = 35.91
$ws.Range("AB5").Value = 1707.82
$ws.Range("AC5").Value = 400
$ws.Range("AD5").Value = 14.74
$ws.Range("AE5").Value = 9685
$ws.Range("AF5").Value = 0.61
$ws.Range("AG5").Value = 180
$ws.Range("AH5").Value = 3.05
$ws.Range("AI5").Value = 33.72
$ws.Range("AJ5").Value = 10000000

# Row 6: update values
$ws.Range("D6").Value = 822
$ws.Range("E6").Value = 56
$ws.Range("F6").Value = 56
$ws.Range("G6").Value = 69
$ws.Range("H6").Value = 54
$ws.Range("I6").Value = 54
$ws.Range("K6").Value = 984
$ws.Range("L6").Value = 215
$ws.Range("M6").Value = 768
$ws.Range("N6").Value = 768
$ws.Range("P6").Value = 50
$ws.Range("Q6").Value = 57
$ws.Range("R6").Value = 45
$ws.Range("S6").Value = -57
$ws.Range("T6").Value = 41
$ws.Range("U6").Value = 16
$ws.Range("V6").Value = 91
$ws.Range("W6").Value = 6.83
$ws.Range("X6").Value = 6.51
$ws.Range("Y6").Value = 7.2
$ws.Range("Z6").Value = 5.42
$ws.Range("AA6").Value = 28.02
$ws.Range("AB6").Value = 1792.25
$ws.Range("AC6").Value = 538
$ws.Range("AD6").Value = 11.55
$ws.Range("AE6").Value = 10241
$ws.Range("AF6").Value = 0.61
$ws.Range("AG6").Value = 200
$ws.Range("AH6").Value = 3.22
$ws.Range("AI6").Value = 27.9
$ws.Range("AJ6").Value = 10000000

# Row 7: clear forecast data (columns D:AI)
$ws.Range("D7:AI7").ClearContents()

# Row 8: clear forecast data (columns D:AI)
$ws.Range("D8:AI8").ClearContents()

# Row 9: clear forecast data (columns D:AI)
$ws.Range("D9:AI9").ClearContents()
